$wb = $excel.ActiveWorkbook

# --- Sheet "Rentometer" ---
$wsRent = $wb.Worksheets.Item("Rentometer")
$wsRent.Range("B17").Value = "https://www.rentometer.com/analysis/3-bed/317-newell-st-barberton-oh-44203/ts4SSOAOOsk/quickview"
$wsRent.Range("B18").Value = 956
$wsRent.Range("B19").Value = "ts4SSOAOOsk"
$wsRent.Range("B20").Value = "[{'rel': 'request pro report', 'href': 'https://www.rentometer.com/api/v1/request_pro_report?api_key=fHSGZM7POi6V5ZPR0w4CXA&token=ts4SSOAOOsk'}, {'rel': 'nearby comps', 'href': 'https://www.rentometer.com/api/v1/nearby_comps?api_key=fHSGZM7POi6V5ZPR0w4CXA&token=ts4SSOAOOsk'}]"

# --- Sheet "Zillow" ---
$wsZillow = $wb.Worksheets.Item("Zillow")
$wsZillow.Range("B3").Value = 1452
$wsZillow.Range("B4").Value = 879
$wsZillow.Range("B5").Value = 1538
$wsZillow.Range("B8").Value = 1319.5
$wsZillow.Range("B9").Value = 1429.5
$wsZillow.Range("B10").Value = 1369

# --- Sheet "rentometer_zillow_user_avg_est" ---
$wsAvg = $wb.Worksheets.Item("rentometer_zillow_user_avg_est")
$wsAvg.Range("B1").Value = 1371.666666666667
$wsAvg.Range("B2").Value = 1365.666666666667
$wsAvg.Range("B3").Value = 1261.75
$wsAvg.Range("B4").Value = 1358.25
